# PAS-6576: Update "individual VIN retrieval" logic to use ENTRY DATE and VALID.
# Adjust the VIN upload test fixture accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: bump the ENTRY DATE (AI) from 2000-01-01 to 2001-01-01.
$ws.Range("AI2").Value = 20010101

# Rows 3-5: replace the old placeholder "Gt" text in column F (VALID-related
# column) with new distinguishable validity markers used by the updated
# "individual VIN retrieval" test logic.
$ws.Range("F3").Value = "invalidVin"
$ws.Range("F4").Value = "SecondValid"
$ws.Range("F5").Value = "ThirdValid"

# Give column F (now holding longer text) an explicit custom width.
$ws.Columns.Item(6).ColumnWidth = 10.3

# The sheet view no longer scrolls the frozen/top-left cell to column S, and
# the remembered selection moves from Z5 to M15.
$ws.Range("M15").Select()
